$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.6
$ws.Range("I2").Value = 2.2
$ws.Range("K2").Value = 1.91
$ws.Range("U2").Value = 2.05
$ws.Range("V2").Value = 1.7
$ws.Range("AC2").Value = 6.5
$ws.Range("AH2").Value = 9.5
$ws.Range("AI2").Value = 10
$ws.Range("AJ2").Value = 21
$ws.Range("AN2").Value = 5
$ws.Range("G4").Value = 1.53
$ws.Range("L4").Value = 7
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 8.5
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.65
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.5
$ws.Range("Z4").Value = 10
$ws.Range("AB4").Value = 41
$ws.Range("AE4").Value = 23
$ws.Range("AH4").Value = 29
$ws.Range("AO4").Value = 8
$ws.Range("AT4").Value = 2.5
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 201
$ws.Range("J5").Value = 8
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("Y5").Value = 26
$ws.Range("AD5").Value = 8
$ws.Range("AJ5").Value = 9
$ws.Range("G6").Value = 1.25
$ws.Range("H6").Value = 5.25
$ws.Range("I6").Value = 13
$ws.Range("L6").Value = 11
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("U6").Value = 2.5
$ws.Range("V6").Value = 1.5
$ws.Range("Z6").Value = 7
$ws.Range("AE6").Value = 29
$ws.Range("AF6").Value = 101
$ws.Range("AG6").Value = 23
$ws.Range("AK6").Value = 101
$ws.Range("AN6").Value = 3
$ws.Range("AO6").Value = 5.5
$ws.Range("AZ6").Value = 351
$ws.Range("BA6").Value = 351
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("Q9").Value = 2.1
$ws.Range("R9").Value = 1.7
$ws.Range("G15").Value = 4.33
$ws.Range("I15").Value = 1.9
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 1.91
$ws.Range("L15").Value = 2.75
$ws.Range("Q15").Value = 2.6
$ws.Range("R15").Value = 1.48
$ws.Range("U15").Value = 2.2
$ws.Range("V15").Value = 1.62
$ws.Range("W15").Value = 9.5
$ws.Range("X15").Value = 21
$ws.Range("Y15").Value = 17
$ws.Range("Z15").Value = 51
$ws.Range("AE15").Value = 21
$ws.Range("AF15").Value = 81
$ws.Range("AH15").Value = 7.5
$ws.Range("AK15").Value = 19
$ws.Range("AN15").Value = 6
$ws.Range("AO15").Value = 29
$ws.Range("AQ15").Value = 101
$ws.Range("AR15").Value = 151
$ws.Range("AS15").Value = 451
$ws.Range("AX15").Value = 11
$ws.Range("BD15").Value = 126
$ws.Range("I16").Value = 2.4
$ws.Range("J16").Value = 3.4
$ws.Range("AA16").Value = 21
$ws.Range("AK16").Value = 21
